$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AB2:AB51").Formula = "=1.1 * AA2"
$ws.Range("AB53:AB81").Formula = "=1.1 * AA53"
